$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns to clean/English machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "de"/"del"/"la"/"los"/"el" -> capitalized "De"/"Del"/"La"/"Los"/"El" in
# municipality / state names (title-casing of connector words)
$ws.Range("B13").Value = "Hidalgo Del Parral"
$ws.Range("A18").Value = "Ciudad De México"
$ws.Range("A27").Value = "Estado De México"
$ws.Range("B28").Value = "Ixtapan De La Sal"
$ws.Range("B30").Value = "Tlalnepantla De Baz"
$ws.Range("B34").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B40").Value = "Ajuchitlán Del Progreso"
$ws.Range("B41").Value = "Atoyac De Álvarez"
$ws.Range("B42").Value = "Chilpancingo De Los Bravo"
$ws.Range("B44").Value = "Técpan De Galeana"
$ws.Range("B48").Value = "Pachuca De Soto"
$ws.Range("B56").Value = "San Miguel El Alto"
$ws.Range("B58").Value = "Zapotitlán De Vadillo"
$ws.Range("B72").Value = "Tepexi De Rodríguez"
$ws.Range("B75").Value = "Pinal De Amoles"

# Minor floating point recompute of percentage totals
$ws.Range("D39").Value = 0.09523809523809525
$ws.Range("D103").Value = 0.09523809523809525

# Remove trailing metadata/footer rows that don't belong to the data table
# (delete higher row numbers first so row indices of the other block don't shift)
$ws.Range("A476:D480").EntireRow.Delete()
$ws.Range("A106:D110").EntireRow.Delete()
